# Auto-generated Excel COM-interop script applying the Ifrit_Profits crafting-profit update.
# For each affected leve row (across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets),
# refresh the market-price-derived columns (H/I/J/K/L/M/N) to the new scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3272.4138
$ws.Range("I76").Value = 3272.4138
$ws.Range("K76").Value = 3272.4138
$ws.Range("M76").Value = -2957.4138

$ws.Range("H79").Value = 3272.4138
$ws.Range("I79").Value = 3272.4138
$ws.Range("K79").Value = 3272.4138
$ws.Range("M79").Value = -2180.4138

$ws.Range("H116").Value = 3090.9092
$ws.Range("I116").Value = 2400
$ws.Range("J116").Value = 3294.1177
$ws.Range("K116").Value = 2400
$ws.Range("L116").Value = 3294.1177
$ws.Range("M116").Value = 1042
$ws.Range("N116").Value = -10178.1177

$ws.Range("H132").Value = 325929.25
$ws.Range("I132").Value = 388223.34
$ws.Range("K132").Value = 1164670.02
$ws.Range("M132").Value = -1162140.02

$ws.Range("H138").Value = 2376.75
$ws.Range("I138").Value = 2014.875
$ws.Range("J138").Value = 2686.9285
$ws.Range("K138").Value = 6044.625
$ws.Range("L138").Value = 8060.7855
$ws.Range("M138").Value = -904.625
$ws.Range("N138").Value = -18340.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 38464190
$ws.Range("I3").Value = 166668160
$ws.Range("K3").Value = 166668160
$ws.Range("M3").Value = -166668045

$ws.Range("H32").Value = 4676.3936
$ws.Range("I32").Value = 5404.5117
$ws.Range("K32").Value = 5404.5117
$ws.Range("M32").Value = -5117.5117

$ws.Range("H45").Value = 72334.92999999999
$ws.Range("I45").Value = 143586
$ws.Range("J45").Value = 1083.8572
$ws.Range("K45").Value = 143586
$ws.Range("L45").Value = 1083.8572
$ws.Range("M45").Value = -143209
$ws.Range("N45").Value = -1837.8572

$ws.Range("H74").Value = 3482.7778
$ws.Range("I74").Value = 934.75
$ws.Range("J74").Value = 4409.3335
$ws.Range("K74").Value = 934.75
$ws.Range("L74").Value = 4409.3335
$ws.Range("M74").Value = -60.75
$ws.Range("N74").Value = -6157.3335

$ws.Range("H77").Value = 3482.7778
$ws.Range("I77").Value = 934.75
$ws.Range("J77").Value = 4409.3335
$ws.Range("K77").Value = 4673.75
$ws.Range("L77").Value = 22046.6675
$ws.Range("M77").Value = -305.75
$ws.Range("N77").Value = -30782.6675

$ws.Range("H103").Value = 42325.8
$ws.Range("J103").Value = 42325.8
$ws.Range("L103").Value = 42325.8
$ws.Range("N103").Value = -44669.8

$ws.Range("H122").Value = 1547.7333
$ws.Range("I122").Value = 1619.5454
$ws.Range("K122").Value = 4858.6362
$ws.Range("M122").Value = -2408.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -87

$ws.Range("H7").Value = 16669251
$ws.Range("I7").Value = 100000000
$ws.Range("J7").Value = 3100.8
$ws.Range("K7").Value = 100000000
$ws.Range("L7").Value = 3100.8
$ws.Range("M7").Value = -99999887
$ws.Range("N7").Value = -3326.8

$ws.Range("H103").Value = 28999.5
$ws.Range("J103").Value = 28999.5
$ws.Range("L103").Value = 28999.5
$ws.Range("N103").Value = -31343.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1550
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = -3226

$ws.Range("H3").Value = 9800
$ws.Range("J3").Value = 9800
$ws.Range("L3").Value = 9800
$ws.Range("N3").Value = -10026

$ws.Range("H4").Value = 7988
$ws.Range("J4").Value = 7988
$ws.Range("L4").Value = 7988
$ws.Range("N4").Value = -8212

$ws.Range("H6").Value = 1674333.6
$ws.Range("I6").Value = 2008600.4
$ws.Range("K6").Value = 2008600.4
$ws.Range("M6").Value = -2008487.4

$ws.Range("H35").Value = 1054.9
$ws.Range("I35").Value = 1054.9
$ws.Range("K35").Value = 1054.9
$ws.Range("M35").Value = -760.9000000000001

$ws.Range("H95").Value = 14200
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 14200
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 14200
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -19692

$ws.Range("H96").Value = 18428.285
$ws.Range("J96").Value = 18428.285
$ws.Range("L96").Value = 18428.285
$ws.Range("N96").Value = -23920.285

$ws.Range("H106").Value = 6000
$ws.Range("J106").Value = 6000
$ws.Range("L106").Value = 6000
$ws.Range("N106").Value = -8524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1757134.4
$ws.Range("J131").Value = 1962644.8
$ws.Range("L131").Value = 5887934.4
$ws.Range("N131").Value = -5898014.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 1750
$ws.Range("J25").Value = 1750
$ws.Range("L25").Value = 1750
$ws.Range("N25").Value = -2808

$ws.Range("H54").Value = 16500
$ws.Range("I54").Value = 3000
$ws.Range("J54").Value = 30000
$ws.Range("K54").Value = 3000
$ws.Range("L54").Value = 30000
$ws.Range("M54").Value = -2610
$ws.Range("N54").Value = -30780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7754000.5
$ws.Range("J2").Value = 10005334
$ws.Range("L2").Value = 10005334
$ws.Range("N2").Value = -10005558

$ws.Range("H62").Value = 22299.666
$ws.Range("J62").Value = 22299.666
$ws.Range("L62").Value = 22299.666
$ws.Range("N62").Value = -23547.666

$ws.Range("H64").Value = 20716.666
$ws.Range("J64").Value = 20716.666
$ws.Range("L64").Value = 20716.666
$ws.Range("N64").Value = -21166.666

$ws.Range("H65").Value = 22299.666
$ws.Range("J65").Value = 22299.666
$ws.Range("L65").Value = 66898.99800000001
$ws.Range("N65").Value = -73138.99800000001

$ws.Range("H67").Value = 20716.666
$ws.Range("J67").Value = 20716.666
$ws.Range("L67").Value = 20716.666
$ws.Range("N67").Value = -22276.666

$ws.Range("H82").Value = 1978.2858
$ws.Range("I82").Value = 1762.5
$ws.Range("J82").Value = 2266
$ws.Range("K82").Value = 1762.5
$ws.Range("L82").Value = 2266
$ws.Range("M82").Value = -1401.5
$ws.Range("N82").Value = -2988

$ws.Range("H85").Value = 1978.2858
$ws.Range("I85").Value = 1762.5
$ws.Range("J85").Value = 2266
$ws.Range("K85").Value = 1762.5
$ws.Range("L85").Value = 2266
$ws.Range("M85").Value = -514.5
$ws.Range("N85").Value = -4762

$ws.Range("H122").Value = 7397.143
$ws.Range("J122").Value = 2775
$ws.Range("L122").Value = 8325
$ws.Range("N122").Value = -13225

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7524.5293
$ws.Range("I132").Value = 11867.111
$ws.Range("J132").Value = 2639.125
$ws.Range("K132").Value = 35601.333
$ws.Range("L132").Value = 7917.375
$ws.Range("M132").Value = -33071.333
$ws.Range("N132").Value = -12977.375

$ws.Range("H138").Value = 39679
$ws.Range("J138").Value = 39679
$ws.Range("L138").Value = 39679
$ws.Range("N138").Value = -49959
